$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Hallo"
$ws.Range("C8").Select()
